$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text-preserving format first, then assign literal values,
# so numeric-looking strings (prices, percentages, counters) are not
# auto-converted to numbers/percentages by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '243.56'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.13%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '5'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '29.90'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '13.66%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '5'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.119'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.41%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '5'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.51%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '5'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.83%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '5'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8393'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.48%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '5'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8639'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '3.28%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '5'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'One'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.01010'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1,593.87%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '5'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1335'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.32%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '5'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06911'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.07%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '5'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02863'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.81%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '5'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09372'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.06%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '5'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001510'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.08%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '5'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04161'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-10.76%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '5'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006009'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.90%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '5'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.506'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-4.05%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '5'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.023'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.29%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '5'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.241'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.65%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '5'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '5'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.03275'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '5.93%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '5'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1296'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '5'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.638'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-2.90%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '5'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.08%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '5'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001209'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-3.10%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '5'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004315'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-4.06%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '5'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001179'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '22.84%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '5'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '0.27%'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '5'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '5'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '5'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '5'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '5'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '5'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '5'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '5'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '5'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '5'
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '5'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '5'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03716'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.12%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '5'
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005326'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-13.81%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '5'
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1057'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.55%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '5'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002086'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-13.07%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '5'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009705'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '9.47%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '5'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005077'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-4.78%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '5'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.05%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '5'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-30.59%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '5'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002722'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '16.57%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '5'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.05%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '5'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.05%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '5'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '5'
